$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The library preparer ("Retrofitted_3569" placeholder) is being replaced
# throughout the data rows (2-21): column B (libraryPreparer) gets the
# actual preparer's initials/name, and column E (purpose) gets the real
# purpose of the run, replacing the same placeholder text that had been
# used for both columns.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 2).Value = "H.BROWN"
    $ws.Cells.Item($r, 5).Value = "fullRNASEQ"
}

# Reflect the author's final on-screen selection/scroll position from the
# edit session.
$ws.Range("D22:F29").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 18
